$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 391
$ws1.Range("F3").Value = 2218
$ws1.Range("F4").Value = 111

# Sheet "全部类型" (sheet4): update the same rows (same events repeated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 391
$ws4.Range("F7").Value = 2218
$ws4.Range("F8").Value = 111
